$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = "DFO_IC_EXISTING"
    3  = "NUC_ST_EXISTING"
    4  = "NG_GT_EXISTING"
    5  = "NG_CC_EXISTING"
    6  = "BLQ_ST_EXISTING"
    7  = "SUN_PV_EXISTING"
    8  = "MWH_BA1H_EXISTING"
    9  = "DFO_GT_EXISTING"
    10 = "WDS_ST_EXISTING"
    11 = "WH_ST_EXISTING"
    12 = "LFG_IC_EXISTING"
    13 = "WND_WT_EXISTING"
    14 = "AB_ST_EXISTING"
    15 = "NG_ST_EXISTING"
    16 = "WAT_HY_EXISTING"
    17 = "WAT_PS_EXISTING"
    18 = "DFO_CC_EXISTING"
    19 = "BIT_ST_EXISTING"
    20 = "LFG_GT_EXISTING"
    21 = "OBG_IC_EXISTING"
    22 = "MWH_BA2H_EXISTING"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
